$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.685507
$ws.Range("H2").Value = 11.056521
$ws.Range("I2").Value = 0.3585631737883472
$ws.Range("J2").Value = 0.3585631737883472
$ws.Range("M2").Value = 0.09551033333333332
$ws.Range("N2").Value = 0.286531
$ws.Range("O2").Value = 0.0198020999427218
$ws.Range("P2").Value = 0.0198020999427218
$ws.Range("Q2").Value = 0.3520040020723333
$ws.Range("R2").Value = 3.168036018651
$ws.Range("S2").Value = 0.007100303803136376
$ws.Range("T2").Value = 0.007100303803136376

# Row 3
$ws.Range("G3").Value = 3.685507
$ws.Range("H3").Value = 11.056521
$ws.Range("I3").Value = 0.3585631737883472
$ws.Range("J3").Value = 0.3585631737883472
$ws.Range("O3").Value = 0.07175622098770619
$ws.Range("P3").Value = 0.07175622098770619
$ws.Range("Q3").Value = 1.275545373183667
$ws.Range("R3").Value = 11.479908358653
$ws.Range("S3").Value = 0.02572913833640994
$ws.Range("T3").Value = 0.02572913833640994

# Row 4
$ws.Range("G4").Value = 3.685507
$ws.Range("H4").Value = 11.056521
$ws.Range("I4").Value = 0.3585631737883472
$ws.Range("J4").Value = 0.3585631737883472
$ws.Range("M4").Value = 4.381634666666667
$ws.Range("N4").Value = 13.144904
$ws.Range("O4").Value = 0.9084416790695721
$ws.Range("P4").Value = 0.9084416790695721
$ws.Range("Q4").Value = 16.14854523544266
$ws.Range("R4").Value = 145.336907118984
$ws.Range("S4").Value = 0.3257337316488009
$ws.Range("T4").Value = 0.3257337316488009

# Row 5
$ws.Range("I5").Value = 0.00964718443071163
$ws.Range("J5").Value = 0.00964718443071163
$ws.Range("M5").Value = 0.09551033333333332
$ws.Range("N5").Value = 0.286531
$ws.Range("O5").Value = 0.0198020999427218
$ws.Range("P5").Value = 0.0198020999427218
$ws.Range("Q5").Value = 0.009470709142999998
$ws.Range("R5").Value = 0.08523638228699999
$ws.Range("S5").Value = 0.0001910345102628214
$ws.Range("T5").Value = 0.0001910345102628214

# Row 6
$ws.Range("I6").Value = 0.00964718443071163
$ws.Range("J6").Value = 0.00964718443071163
$ws.Range("O6").Value = 0.07175622098770619
$ws.Range("P6").Value = 0.07175622098770619
$ws.Range("S6").Value = 0.0006922454979193022
$ws.Range("T6").Value = 0.0006922454979193022

# Row 7
$ws.Range("I7").Value = 0.00964718443071163
$ws.Range("J7").Value = 0.00964718443071163
$ws.Range("M7").Value = 4.381634666666667
$ws.Range("N7").Value = 13.144904
$ws.Range("O7").Value = 0.9084416790695721
$ws.Range("P7").Value = 0.9084416790695721
$ws.Range("Q7").Value = 0.434478511912
$ws.Range("R7").Value = 3.910306607208
$ws.Range("S7").Value = 0.008763904422529506
$ws.Range("T7").Value = 0.008763904422529506

# Row 8
$ws.Range("G8").Value = 6.493877
$ws.Range("H8").Value = 19.481631
$ws.Range("I8").Value = 0.6317896417809412
$ws.Range("J8").Value = 0.6317896417809411
$ws.Range("M8").Value = 0.09551033333333332
$ws.Range("N8").Value = 0.286531
$ws.Range("O8").Value = 0.0198020999427218
$ws.Range("P8").Value = 0.0198020999427218
$ws.Range("Q8").Value = 0.6202323568956666
$ws.Range("R8").Value = 5.582091212061
$ws.Range("S8").Value = 0.0125107616293226
$ws.Range("T8").Value = 0.0125107616293226

# Row 9
$ws.Range("G9").Value = 6.493877
$ws.Range("H9").Value = 19.481631
$ws.Range("I9").Value = 0.6317896417809412
$ws.Range("J9").Value = 0.6317896417809411
$ws.Range("O9").Value = 0.07175622098770619
$ws.Range("P9").Value = 0.07175622098770619
$ws.Range("Q9").Value = 2.247515677320333
$ws.Range("R9").Value = 20.227641095883
$ws.Range("S9").Value = 0.04533483715337695
$ws.Range("T9").Value = 0.04533483715337694

# Row 10
$ws.Range("G10").Value = 6.493877
$ws.Range("H10").Value = 19.481631
$ws.Range("I10").Value = 0.6317896417809412
$ws.Range("J10").Value = 0.6317896417809411
$ws.Range("M10").Value = 4.381634666666667
$ws.Range("N10").Value = 13.144904
$ws.Range("O10").Value = 0.9084416790695721
$ws.Range("P10").Value = 0.9084416790695721
$ws.Range("Q10").Value = 28.45379658426933
$ws.Range("R10").Value = 256.084169258424
$ws.Range("S10").Value = 0.5739440429982418
$ws.Range("T10").Value = 0.5739440429982416
